# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever the value is currently "System, dnasr281@gmail.com".
# Other variants (e.g. already "dnasr281@gmail.com, System" or just the
# bare email address) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
